# Update the "Players" sheet: game clock status strings and updated box-score stats
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

$ws.Range("G3").Value = "1:01 - 2nd Half"
$ws.Range("O3").Value = 33

$ws.Range("G4").Value = "1:01 - 2nd Half"

$ws.Range("O5").Value = 16

$ws.Range("G6").Value = "1:01 - 2nd Half"
$ws.Range("H6").Value = 14
$ws.Range("M6").Value = 2
$ws.Range("O6").Value = 35

$ws.Range("G7").Value = "1:01 - 2nd Half"
$ws.Range("O7").Value = 32

$ws.Range("G8").Value = "1:01 - 2nd Half"
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 4
$ws.Range("O8").Value = 32

$ws.Range("G11").Value = "1:01 - 2nd Half"

$ws.Range("O13").Value = 17

$ws.Range("G14").Value = "1:01 - 2nd Half"

$ws.Range("G15").Value = "1:01 - 2nd Half"

$ws.Range("G18").Value = "1:01 - 2nd Half"
$ws.Range("H18").Value = 23
$ws.Range("I18").Value = 21
$ws.Range("J18").Value = 3
$ws.Range("O18").Value = 36

$ws.Range("G19").Value = "1:01 - 2nd Half"
$ws.Range("H19").Value = 12
$ws.Range("O19").Value = 38

$ws.Range("G20").Value = "1:01 - 2nd Half"

$ws.Range("G21").Value = "1:01 - 2nd Half"
$ws.Range("H21").Value = 18
$ws.Range("J21").Value = 12
$ws.Range("O21").Value = 28

$ws.Range("G23").Value = "1:01 - 2nd Half"
$ws.Range("O23").Value = 35

$ws.Range("H24").Value = 8
$ws.Range("I24").Value = 9

$ws.Range("G25").Value = "1:01 - 2nd Half"
$ws.Range("O25").Value = 26

$ws.Range("G26").Value = "1:01 - 2nd Half"
$ws.Range("O26").Value = 19

$ws.Range("O27").Value = 11

$ws.Range("G33").Value = "1:01 - 2nd Half"

# Update the "OwnerTotals" sheet: starter pooh totals + owner ordering for tied ranks
$ws2 = $wb.Worksheets.Item("OwnerTotals")

$ws2.Range("B2").Value = 23

$ws2.Range("B3").Value = 18

$ws2.Range("A4").Value = "CDL"
$ws2.Range("B4").Value = 14

$ws2.Range("A5").Value = "Booz"

